$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.345.43"
$ws.Range("E2").Value = "  +3.85%  "
$ws.Range("D3").Value = "3.490.32"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Formula = "'579.77"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Formula = "'160.88"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Formula = "'0.606"
$ws.Range("E8").Value = "  +11.88%  "
$ws.Range("D9").Value = "3.491.51"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("E11").Value = "  +3.96%  "
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("D13").Value = "4.093.21"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("E15").Value = "  +3.59%  "
$ws.Range("D16").Formula = "'28.71"
$ws.Range("E16").Value = "  +7.05%  "
$ws.Range("D17").Value = "65.388.15"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "3.486.87"
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("E19").Value = "  +4.31%  "
$ws.Range("D20").Formula = "'14.40"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("D21").Formula = "'387.18"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").Formula = "'8.29"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("D23").Formula = "'0.555"
$ws.Range("E23").Value = "  +5.31%  "
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").Formula = "'0.999"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +6.34%  "
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Formula = "'1.49"
$ws.Range("E30").Value = "  +12.46%  "
$ws.Range("E31").Value = "  +4.92%  "
$ws.Range("E32").Value = "  +4.05%  "
$ws.Range("D33").Formula = "'23.68"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Formula = "'7.29"
$ws.Range("E34").Value = "  +8.34%  "
$ws.Range("E35").Value = "  +9.71%  "
$ws.Range("D36").Formula = "'162.87"
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("E37").Value = "  +6.21%  "
$ws.Range("D38").Value = "3.000.78"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Formula = "'0.0780"
$ws.Range("E39").Value = "  +3.19%  "
$ws.Range("D40").Formula = "'27.26"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").Formula = "'6.58"
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("D44").Formula = "'42.81"
$ws.Range("E44").Value = "  +4.27%  "
$ws.Range("D45").Formula = "'0.780"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").Formula = "'25.65"
$ws.Range("E46").Value = "  +10.78%  "
$ws.Range("D47").Formula = "'1.12"
$ws.Range("E47").Value = "  +5.28%  "
$ws.Range("D48").Formula = "'324.50"
$ws.Range("E48").Value = "  +13.77%  "
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("D50").Formula = "'2.22"
$ws.Range("E50").Value = "  +4.42%  "
$ws.Range("E51").Value = "  +6.73%  "
